$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

Set-TextValue 2 4 '67.940.96'
Set-TextValue 2 5 '  +0.78%  '
Set-TextValue 3 4 '2.639.68'
Set-TextValue 3 5 '  +0.55%  '
Set-TextValue 4 5 '  +0.01%  '
Set-TextValue 5 4 '598.29'
Set-TextValue 5 5 '  -0.06%  '
Set-TextValue 6 4 '153.77'
Set-TextValue 6 5 '  +0.48%  '
Set-TextValue 7 5 '  +0.01%  '
Set-TextValue 8 4 '0.551'
Set-TextValue 8 5 '  -0.46%  '
Set-TextValue 9 4 '2.638.85'
Set-TextValue 9 5 '  +0.56%  '
Set-TextValue 10 5 '  +10.36%  '
Set-TextValue 11 5 '  -0.70%  '
Set-TextValue 12 5 '  +0.62%  '
Set-TextValue 13 4 '0.348'
Set-TextValue 13 5 '  -0.19%  '
Set-TextValue 14 4 '27.69'
Set-TextValue 14 5 '  +0.14%  '
Set-TextValue 15 5 '  +3.39%  '
Set-TextValue 16 4 '3.119.96'
Set-TextValue 16 5 '  +0.85%  '
Set-TextValue 17 4 '67.835.21'
Set-TextValue 17 5 '  +0.65%  '
Set-TextValue 18 4 '2.633.49'
Set-TextValue 18 5 '  +0.74%  '
Set-TextValue 19 4 '11.46'
Set-TextValue 19 5 '  +2.83%  '
Set-TextValue 20 4 '373.56'
Set-TextValue 20 5 '  +2.86%  '
Set-TextValue 21 5 '  +0.04%  '
Set-TextValue 22 5 '  -0.94%  '
Set-TextValue 23 5 '  -1.60%  '
Set-TextValue 24 5 '  -1.93%  '
Set-TextValue 25 4 '72.23'
Set-TextValue 25 5 '  +1.64%  '
Set-TextValue 26 5 '  -0.04%  '
Set-TextValue 27 4 '9.95'
Set-TextValue 27 5 '  -1.26%  '
Set-TextValue 28 4 '2.760.97'
Set-TextValue 28 5 '  -0.01%  '
Set-TextValue 29 5 '  +1.95%  '
Set-TextValue 30 5 '  -0.42%  '
Set-TextValue 31 4 '579.97'
Set-TextValue 31 5 '  -0.77%  '
Set-TextValue 32 5 '  +0.22%  '
Set-TextValue 33 5 '  +0.58%  '
Set-TextValue 34 5 '  +0.42%  '
Set-TextValue 35 5 '  +0.06%  '
Set-TextValue 36 5 '  +0.15%  '
Set-TextValue 37 4 '1.52'
Set-TextValue 37 5 '  -0.40%  '
Set-TextValue 38 4 '157.94'
Set-TextValue 38 5 '  +0.45%  '
Set-TextValue 39 4 '19.21'
Set-TextValue 39 5 '  +0.25%  '
Set-TextValue 40 4 '1.90'
Set-TextValue 40 5 '  +4.98%  '
Set-TextValue 41 5 '  +0.16%  '
Set-TextValue 42 4 '5.37'
Set-TextValue 42 5 '  +1.69%  '
Set-TextValue 43 4 '0.0₆0343'
Set-TextValue 43 5 '  +19.31%  '
Set-TextValue 44 5 '  +2.61%  '
Set-TextValue 45 4 '17.15'
Set-TextValue 45 5 '  +4.90%  '
Set-TextValue 46 5 '  +0.06%  '
Set-TextValue 47 4 '40.24'
Set-TextValue 47 5 '  -2.20%  '
Set-TextValue 48 4 '156.32'
Set-TextValue 48 5 '  -0.16%  '
Set-TextValue 49 4 '3.70'
Set-TextValue 49 5 '  -0.94%  '
Set-TextValue 50 4 '21.97'
Set-TextValue 50 5 '  +0.97%  '
Set-TextValue 51 4 '1.71'
Set-TextValue 51 5 '  -1.54%  '
